$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily data point for 2026/01/07 (hour 20) was appended to the log.
# It sorts right after the existing 2026/01/07 rows (596-600) and before the
# 2026/12/29 block, so insert a fresh row at 601, pushing rows 601:642 down
# to 602:643.
$ws.Rows.Item(601).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds a date-looking string stored as literal text (matches the
# rest of the sheet, which uses text dates rather than real date values), so
# force text interpretation with a leading apostrophe and strip the
# auto-applied text style back to Normal so no stray style id is introduced.
$ws.Cells.Item(601, 1).Value = "'2026/01/07"
$ws.Cells.Item(601, 1).Style = "Normal"
$ws.Cells.Item(601, 2).Value = "水"
$ws.Cells.Item(601, 3).Value = 20
$ws.Cells.Item(601, 4).Value = 201
